$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Row 31: new commit entry "Smaller changes" -> compile.bat
# ---------------------------------------------------------------------------
$ws.Range("A31").Value = "Smaller changes"
$ws.Range("B31").Value = "compile.bat"
$ws.Range("C31").Value = "everything"
$ws.Range("D31").Value = 9
$ws.Range("E31").Value = 2
$ws.Range("F31").Value = "none"
$ws.Range("G31").Value = "adalfarus"

# Turn A31 into a commit hyperlink (same look as the other commit cells)
$ws.Hyperlinks.Add($ws.Range("A31"), "https://github.com/Giesbrt/Automaten/commit/8f1a3b6c9d2e4f5061728394a5b6c7d8e9f0a1b2") | Out-Null
# Re-apply the exact "commit link" formatting used by the sibling rows
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A31").Value = "Smaller changes"

# ---------------------------------------------------------------------------
# Row 32: continuation of the same commit -> pyautoinst-config.json
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = ".."
$ws.Range("B32").Value = "pyautoinst-config.json"
$ws.Range("C32").Value = "everything"
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = "none"
$ws.Range("G32").Value = "adalfarus"
$ws.Range("H32").Value = "Had some edge cases"

# ---------------------------------------------------------------------------
# Sheet view: scroll position / selection
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E31").Select() | Out-Null
